$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Date: 11/28/2021 -> 12/04/2021
# ---------------------------------------------------------------------
$d.Content.Find.Execute("11/28/2021", $false, $false, $false, $false, $false, `
    $true, 1, $false, "12/04/2021", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. modelling -> modeling
# ---------------------------------------------------------------------
$d.Content.Find.Execute("statistical modelling requires complete data", $false, $false, $false, $false, $false, `
    $true, 1, $false, "statistical modeling requires complete data", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. "Base R: plot(), pairs()" -> "Base R: Boxplots(), pairs()" and
#    remove the trailing manual line break run that followed it.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Base R: plot(), pairs()", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Base R: Boxplots(), pairs()", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Base R: Boxplots(), pairs()", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$brRange = $d.Range($rng.End, $rng.End + 1)
if ($brRange.Text -eq [char]11) {
    $brRange.Delete()
}

# ---------------------------------------------------------------------
# 4. Remove the trailing manual line break run after the first
#    "'Missing Data' script" occurrence.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$needle = [char]0x2018 + "Missing Data" + [char]0x2019 + " script"
$rng2.Find.Execute($needle, $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$brRange2 = $d.Range($rng2.End, $rng2.End + 1)
if ($brRange2.Text -eq [char]11) {
    $brRange2.Delete()
}

# ---------------------------------------------------------------------
# 5. Boxplots/zRange -> Boxplot/zRange
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Boxplots/zRange/plausible meanSD", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Boxplot/zRange/plausible meanSD", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. "1. Model Assumptions will are not met..." -> "1. Model assumptions are not met..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("1. Model Assumptions will are not met for each IRT model to be fit, including monotonicity, local dependence and item invariance", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "1. Model assumptions are not met for each IRT model to be fit, including monotonicity, local dependence and item invariance", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. Append the new "timeline" content at the end of the document,
#    right after the "Please provide a brief timeline..." heading and
#    before the closing bookmark / sectPr.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$insertionPoint.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last
$p1.Style = "FirstParagraph"
$p1.Range.InsertAfter("Planning and writing for the preregistration started in October, 2021. Data was simulated based on the codebook published on OSF MARK_OSF_LINK in November, 2021 and published on GitHub under a branch of the " + [char]0x201C + "Prereg_Function" + [char]0x201D + " project, called " + [char]0x201C + "ClothingIRT" + [char]0x201D + " - (MARK_GITHUB_LINK). Set-up and hypothesis building will completed in December with results of additional literature review and simulated data/coding embedded in this parent RMarkdown document. The actual research below is planned for completion in December 2021.")

$insertionPoint2 = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$insertionPoint2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = "BodyText"
$p2.Range.InsertAfter("Hypothesis 1 will be tested after completion of the set-up and hypothesis building. The results of Hypothesis 1 will be posted in an expanded version of this RMarkdown document.")

$insertionPoint3 = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$insertionPoint3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = "BodyText"
$p3.Range.InsertAfter("Based on the results of Hypothesis 1 and any requisite additional background research, revisions to Hypothesis 2, coding and results will occur and will be documented in this parent RMarkdown document.")

# Now convert the two placeholder markers into real hyperlinks.
$osfRng = $d.Content
$osfRng.Find.Execute("MARK_OSF_LINK", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($osfRng, "https://osf.io/ajv5z/", $null, $null, "https://osf.io/ajv5z/") | Out-Null

$ghRng = $d.Content
$ghRng.Find.Execute("MARK_GITHUB_LINK", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($ghRng, "https://github.com/zenit125/Prereg_Functions/tree/Clothing_IRT", $null, $null, `
    "https://github.com/zenit125/Prereg_Functions/tree/Clothing_IRT") | Out-Null
